$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 667817.75
$ws.Range("I2").Value = 452.57144
$ws.Range("J2").Value = 1251762.2
$ws.Range("K2").Value = 452.57144
$ws.Range("L2").Value = 1251762.2
$ws.Range("M2").Value = -339.57144
$ws.Range("N2").Value = -1251988.2

$ws.Range("H12").Value = 648.1667
$ws.Range("J12").Value = 662.6667
$ws.Range("L12").Value = 662.6667
$ws.Range("N12").Value = -1002.6667

$ws.Range("H40").Value = 35717380
$ws.Range("J40").Value = 62502840
$ws.Range("L40").Value = 62502840
$ws.Range("N40").Value = -62503190

$ws.Range("H70").Value = 642959.5
$ws.Range("J70").Value = 1389.0834
$ws.Range("L70").Value = 4167.2502
$ws.Range("N70").Value = -4707.2502

$ws.Range("H73").Value = 642959.5
$ws.Range("J73").Value = 1389.0834
$ws.Range("L73").Value = 4167.2502
$ws.Range("N73").Value = -6039.2502

$ws.Range("H129").Value = 3675.7273
$ws.Range("J129").Value = 11660.667
$ws.Range("L129").Value = 34982.001
$ws.Range("N129").Value = -44982.001

$ws.Range("H137").Value = 2169.1052
$ws.Range("I137").Value = 2087.7144
$ws.Range("K137").Value = 6263.1432
$ws.Range("M137").Value = -3713.1432

$ws.Range("H138").Value = 7424.467
$ws.Range("I138").Value = 3811.8
$ws.Range("J138").Value = 14649.8
$ws.Range("K138").Value = 11435.4
$ws.Range("L138").Value = 43949.39999999999
$ws.Range("M138").Value = -6295.400000000001
$ws.Range("N138").Value = -54229.39999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2589
$ws.Range("I2").Value = 2737.1428
$ws.Range("K2").Value = 2737.1428
$ws.Range("M2").Value = -2624.1428

$ws.Range("H61").Value = 20003080
$ws.Range("I61").Value = 26253734
$ws.Range("K61").Value = 26253734
$ws.Range("M61").Value = -26253522

$ws.Range("H63").Value = 4249.3335
$ws.Range("I63").Value = 4339.2
$ws.Range("J63").Value = 3800
$ws.Range("K63").Value = 4339.2
$ws.Range("L63").Value = 3800
$ws.Range("M63").Value = -3653.2
$ws.Range("N63").Value = -5172

$ws.Range("H66").Value = 4249.3335
$ws.Range("I66").Value = 4339.2
$ws.Range("J66").Value = 3800
$ws.Range("K66").Value = 21696
$ws.Range("L66").Value = 19000
$ws.Range("M66").Value = -18264
$ws.Range("N66").Value = -25864

$ws.Range("H92").Value = 66499.5
$ws.Range("J92").Value = 66499.5
$ws.Range("L92").Value = 66499.5
$ws.Range("N92").Value = -71491.5

$ws.Range("H116").Value = 2589
$ws.Range("I116").Value = 2737.1428
$ws.Range("K116").Value = 2737.1428
$ws.Range("M116").Value = -443.1428000000001

$ws.Range("H136").Value = 20003080
$ws.Range("I136").Value = 26253734
$ws.Range("K136").Value = 78761202
$ws.Range("M136").Value = -78758652

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2589
$ws.Range("I3").Value = 2737.1428
$ws.Range("K3").Value = 2737.1428
$ws.Range("M3").Value = -2623.1428

$ws.Range("H20").Value = 12729.625
$ws.Range("I20").Value = 16222.833
$ws.Range("J20").Value = 2250
$ws.Range("K20").Value = 16222.833
$ws.Range("L20").Value = 2250
$ws.Range("M20").Value = -15975.833
$ws.Range("N20").Value = -2744

$ws.Range("H86").Value = 5118.7
$ws.Range("I86").Value = 2228.889
$ws.Range("J86").Value = 7483.091
$ws.Range("K86").Value = 2228.889
$ws.Range("L86").Value = 7483.091
$ws.Range("M86").Value = -1105.889
$ws.Range("N86").Value = -9729.091

$ws.Range("H89").Value = 5118.7
$ws.Range("I89").Value = 2228.889
$ws.Range("J89").Value = 7483.091
$ws.Range("K89").Value = 11144.445
$ws.Range("L89").Value = 37415.455
$ws.Range("M89").Value = -5528.445
$ws.Range("N89").Value = -48647.455

$ws.Range("H105").Value = 740514.0600000001
$ws.Range("I105").Value = 996097.1
$ws.Range("J105").Value = 5712.75
$ws.Range("K105").Value = 996097.1
$ws.Range("L105").Value = 5712.75
$ws.Range("M105").Value = -994350.1
$ws.Range("N105").Value = -9206.75

$ws.Range("H134").Value = 5002926
$ws.Range("I134").Value = 2735.4119
$ws.Range("J134").Value = 33337338
$ws.Range("K134").Value = 8206.235700000001
$ws.Range("L134").Value = 100012014
$ws.Range("M134").Value = -5671.235700000001
$ws.Range("N134").Value = -100017084

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2362.5
$ws.Range("I22").Value = 2250.5
$ws.Range("K22").Value = 2250.5
$ws.Range("M22").Value = -1900.5

$ws.Range("H93").Value = 74089.2
$ws.Range("I93").Value = 70000
$ws.Range("J93").Value = 75111.5
$ws.Range("K93").Value = 70000
$ws.Range("L93").Value = 75111.5
$ws.Range("M93").Value = -68128
$ws.Range("N93").Value = -78855.5

$ws.Range("H107").Value = 1327.2916
$ws.Range("I107").Value = 840.7619
$ws.Range("K107").Value = 840.7619
$ws.Range("M107").Value = 1079.2381

$ws.Range("H132").Value = 2848.6785
$ws.Range("I132").Value = 2642.2
$ws.Range("J132").Value = 3364.875
$ws.Range("K132").Value = 7926.599999999999
$ws.Range("L132").Value = 10094.625
$ws.Range("M132").Value = -5396.599999999999
$ws.Range("N132").Value = -15154.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 5390.4287
$ws.Range("I13").Value = 1033.6666
$ws.Range("K13").Value = 3100.9998
$ws.Range("M13").Value = -2932.9998

$ws.Range("H34").Value = 10918.286
$ws.Range("J34").Value = 15166
$ws.Range("L34").Value = 45498
$ws.Range("N34").Value = -45666

$ws.Range("H92").Value = 142.75
$ws.Range("J92").Value = 142.75
$ws.Range("L92").Value = 428.25
$ws.Range("N92").Value = -2924.25

$ws.Range("H132").Value = 1383.1428
$ws.Range("I132").Value = 855.4
$ws.Range("K132").Value = 7698.599999999999
$ws.Range("M132").Value = -5168.599999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

$ws.Range("H35").Value = 89999
$ws.Range("J35").Value = 89999
$ws.Range("L35").Value = 89999
$ws.Range("N35").Value = -90595

$ws.Range("H70").Value = 10396.875
$ws.Range("I70").Value = 9673.75
$ws.Range("J70").Value = 11120
$ws.Range("K70").Value = 9673.75
$ws.Range("L70").Value = 11120
$ws.Range("M70").Value = -9403.75
$ws.Range("N70").Value = -11660

$ws.Range("H73").Value = 10396.875
$ws.Range("I73").Value = 9673.75
$ws.Range("J73").Value = 11120
$ws.Range("K73").Value = 9673.75
$ws.Range("L73").Value = 11120
$ws.Range("M73").Value = -8737.75
$ws.Range("N73").Value = -12992

$ws.Range("H80").Value = 2907.5715
$ws.Range("J80").Value = 3630.6
$ws.Range("L80").Value = 3630.6
$ws.Range("N80").Value = -5626.6

$ws.Range("H83").Value = 2907.5715
$ws.Range("J83").Value = 3630.6
$ws.Range("L83").Value = 18153
$ws.Range("N83").Value = -28137

$ws.Range("H93").Value = 50000
$ws.Range("J93").Value = 50000
$ws.Range("L93").Value = 50000
$ws.Range("N93").Value = -53744

$ws.Range("H102").Value = 2819.3547
$ws.Range("J102").Value = 2567
$ws.Range("L102").Value = 2567
$ws.Range("N102").Value = -5811

$ws.Range("H132").Value = 7695093
$ws.Range("I132").Value = 2832.1667
$ws.Range("J132").Value = 14288459
$ws.Range("K132").Value = 8496.500100000001
$ws.Range("L132").Value = 42865377
$ws.Range("M132").Value = -5966.500100000001
$ws.Range("N132").Value = -42870437

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5541.6665
$ws.Range("I7").Value = 4768
$ws.Range("J7").Value = 8249.5
$ws.Range("K7").Value = 4768
$ws.Range("L7").Value = 8249.5
$ws.Range("M7").Value = -4656
$ws.Range("N7").Value = -8473.5

$ws.Range("H16").Value = 5671.4287
$ws.Range("I16").Value = 2043
$ws.Range("J16").Value = 14742.5
$ws.Range("K16").Value = 2043
$ws.Range("L16").Value = 14742.5
$ws.Range("M16").Value = -1873
$ws.Range("N16").Value = -15082.5

$ws.Range("H68").Value = 4169138.8
$ws.Range("I68").Value = 8335261
$ws.Range("J68").Value = 3016.4
$ws.Range("K68").Value = 8335261
$ws.Range("L68").Value = 3016.4
$ws.Range("M68").Value = -8334512
$ws.Range("N68").Value = -4514.4

$ws.Range("H71").Value = 4169138.8
$ws.Range("I71").Value = 8335261
$ws.Range("J71").Value = 3016.4
$ws.Range("K71").Value = 41676305
$ws.Range("L71").Value = 15082
$ws.Range("M71").Value = -41672561
$ws.Range("N71").Value = -22570

$ws.Range("H100").Value = 25003704
$ws.Range("J100").Value = 83337330
$ws.Range("L100").Value = 83337330
$ws.Range("N100").Value = -83338412

$ws.Range("H126").Value = 5541.6665
$ws.Range("I126").Value = 4768
$ws.Range("J126").Value = 8249.5
$ws.Range("K126").Value = 14304
$ws.Range("L126").Value = 24748.5
$ws.Range("M126").Value = -11834
$ws.Range("N126").Value = -29688.5

$ws.Range("H136").Value = 3111.5
$ws.Range("I136").Value = 3067.625
$ws.Range("K136").Value = 9202.875
$ws.Range("M136").Value = -6652.875
